$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.509.77"
$ws.Range("E2").Value = "  -0.85%  "
$ws.Range("D3").Value = "3.546.25"
$ws.Range("E3").Value = "  -1.65%  "
$style = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = $style
$ws.Range("E4").Value = "  +0.13%  "
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "197.67"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  +0.84%  "
$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "584.65"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  -3.02%  "
$ws.Range("E7").Value = "  -2.32%  "
$style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = $style
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -1.63%  "
$style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.631"
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = "  -2.69%  "
$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "51.93"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  -3.53%  "
$style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000287"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = "  -5.71%  "
$ws.Range("E13").Value = "  -3.08%  "
$ws.Range("D14").Value = "4.109.52"
$ws.Range("E14").Value = "  -1.83%  "
$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "664.84"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  +12.25%  "
$ws.Range("D16").Value = "69.672.81"
$ws.Range("E16").Value = "  -0.84%  "
$ws.Range("D17").Value = "3.543.24"
$ws.Range("E17").Value = "  -1.67%  "
$ws.Range("E18").Value = "  -5.41%  "
$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.58"
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = "  -3.16%  "
$ws.Range("E20").Value = "  -0.65%  "
$ws.Range("E21").Value = "  -2.51%  "
$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "18.33"
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = "  +3.67%  "
$ws.Range("E23").Value = "  +2.70%  "
$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "105.33"
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = "  +3.18%  "
$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.38"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  -4.65%  "
$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.92"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = "  -3.67%  "
$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.20"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  -5.68%  "
$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.68"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = "  +0.97%  "
$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.49"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  -1.45%  "
$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.47"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  -6.68%  "
$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.85"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  -3.97%  "
$ws.Range("E32").Value = "  -3.42%  "
$ws.Range("E33").Value = "  -4.85%  "
$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "61.97"
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = "  -2.05%  "
$ws.Range("D35").Value = "3.789.69"
$ws.Range("E35").Value = "  -3.46%  "
$ws.Range("B36").Value = "Stacks"
$ws.Range("C36").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.74"
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = "  +5.77%  "
$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").Value = "0.0₃0813"
$ws.Range("E37").Value = "  -8.35%  "
$ws.Range("E38").Value = "  -0.06%  "
$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "505.07"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  -4.62%  "
$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.94"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  -6.06%  "
$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.374"
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = "  -4.53%  "
$ws.Range("E42").Value = "  +0.88%  "
$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "34.72"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = "  -6.47%  "
$ws.Range("E44").Value = "  +0.03%  "
$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.89"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  +1.06%  "
$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.39"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = "  -0.03%  "
$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.136"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  -3.38%  "
$ws.Range("E48").Value = "  -0.17%  "
$ws.Range("E49").Value = "  -3.29%  "
$ws.Range("E50").Value = "  +20.71%  "
$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.68"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  +60.99%  "
